## Job para atualizar pautas
## - Add a collapsed bookmark "__DdeLink__23_1057900407" right before the
##   "[P.POSICIONAMENTO]" run (existing "__DdeLink__115_1156576634" bookmark
##   auto-renumbers from id 0 -> id 1).
## - Insert a new paragraph "[P.EXPLICACAO]" right after the
##   "Posicao do MJ: [P.POSICIONAMENTO]" paragraph (non-bold, ind left=1440).
## - Register three new character styles: "ListLabel 8/9/10".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Bookmark around the insertion point just before [P.POSICIONAMENTO]
# ---------------------------------------------------------------------
$posRng = $d.Content
[void]$posRng.Find.Execute("[P.POSICIONAMENTO]", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posStart = $posRng.Start
$bmRng = $d.Range($posStart, $posStart)
$d.Bookmarks.Add("__DdeLink__23_1057900407", $bmRng)

# ---------------------------------------------------------------------
# 2) New paragraph "[P.EXPLICACAO]" after the "Posicao do MJ" paragraph
# ---------------------------------------------------------------------
$posRng2 = $d.Content
[void]$posRng2.Find.Execute("[P.POSICIONAMENTO]", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mjPara = $posRng2.Paragraphs(1)
$followingPara = $mjPara.Next()

# Inserting before the paragraph that follows "Posicao do MJ: ..." gives a
# new paragraph mark that inherits that following paragraph's (non-listed)
# formatting, so it starts life with no numbering applied.
$followingPara.Range.InsertParagraphBefore()

$newPara = $mjPara.Next()
$newPara.Range.Text = "[P.EXPLICACAO]"
$newPara.Range.ParagraphFormat.LeftIndent = 72
$newPara.Range.Font.Bold = $false
$newPara.Range.Font.BoldBi = $false

# ---------------------------------------------------------------------
# 3) New character styles used by list labels 8-10
# ---------------------------------------------------------------------
$ll8 = $d.Styles.Add("ListLabel 8", 2)
$ll8.Font.NameBi = "OpenSymbol"
$ll8.Font.Underline = 0

$ll9 = $d.Styles.Add("ListLabel 9", 2)
$ll9.Font.NameBi = "Wingdings"
$ll9.Font.Underline = 0

$ll10 = $d.Styles.Add("ListLabel 10", 2)
$ll10.Font.NameBi = "Wingdings 2"
$ll10.Font.Underline = 0
